# Post-summer-2021 update: remove two stale/duplicate "Johan Hjort" (t="s" idx 5)
# toktnr rows (2003703 and 2004703) that were grouping incorrectly, letting the
# remaining rows shift up to close the gap. Also reset the view back to the
# top of the sheet and select the full row that was being edited (row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row currently holding year 2003 / toktnr 2003703 (row 3).
$ws.Rows.Item(3).Delete()

# After that shift, the row that used to hold year 2004 / toktnr 2004703 is
# now row 7 (was row 8) - remove it too.
$ws.Rows.Item(7).Delete()

# Reset the view: scroll back to the top (no frozen/offset topLeftCell) and
# select the full row 7 (A7:XFD7), matching the post-edit selection.
$ws.Rows.Item(7).Select()
